# Update "Kanara Greens and Flowers Trading_2024-9-24.xlsx"
#
# Orders sheet (sheet1):
#   - F2: "18" -> "180" (kept as text, not a number)
#   - New row 3: C3 = "2_粉洋桔梗_Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners" (text)
#   - dimension / used range grows to A1:L3 automatically from the above writes
#
# Summary sheet (sheet2):
#   - G2: "018" -> "01800" (kept as text, not a number)

$wb = $excel.ActiveWorkbook

$orders = $wb.Worksheets.Item("Orders")

# F2 currently holds the text "18" (number-stored-as-text). Force the
# cell format to Text before writing so Excel doesn't coerce the new
# value into a real number (which would also drop any leading zeros).
$orders.Range("F2").NumberFormat = "@"
$orders.Range("F2").Value = "180"

# Add the new second flower row.
$orders.Range("C3").NumberFormat = "@"
$orders.Range("C3").Value = "2_粉洋桔梗_Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners"

$summary = $wb.Worksheets.Item("Summary")

# G2 currently holds the text "018" -> "01800"; keep it textual.
$summary.Range("G2").NumberFormat = "@"
$summary.Range("G2").Value = "01800"
